# Approach description document tidy-up edits.
#
# Before (list items, in order):
#   P15: "Used the frequency reports to detect anomalies ..."
#   P16: "Converted  "unknown" values in the regionType column ... "
#   P17: "In frequency reports, I included the 'missing' option ... newCellUser."
#   P18: "" (final, empty list paragraph)
#
# After:
#   P15: "Used linux command (via git bash) "less" to view data ... import errors."
#   P16: "Used the frequency reports to detect anomalies ..."   (old P15 text)
#   P17: "Converted "unknown" values in the regionType column ... " (old P16 text, minor space fix)
#   P18: "In frequency reports, I included the 'missing' option ... newCellUser." (old P17 text, new paragraph)
#   P19: "" (old P18, but no longer a list item - indented plain paragraph)

$d = $word.ActiveDocument

$lq = [char]0x201C
$rq = [char]0x201D
$lsq = [char]0x2018
$rsq = [char]0x2019

# ---------------------------------------------------------------------------
# Step 1: insert a brand-new (empty) list paragraph right before the
# "In frequency reports..." paragraph (currently paragraph 17). This leaves
# the "In frequency reports..." paragraph itself completely untouched
# (it simply shifts down to become paragraph 18), and gives us a blank
# paragraph 17 to turn into the corrected "Converted ..." text.
# ---------------------------------------------------------------------------
$pInFreq = $d.Paragraphs(17)
$pInFreq.Range.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# Step 2: paragraph 17 is now a brand-new empty paragraph - fill it with the
# corrected version of the old "Converted ..." text (double spaces tidied to
# single spaces), including its own trailing line break.
# ---------------------------------------------------------------------------
$p17 = $d.Paragraphs(17)
$r17 = $p17.Range
$r17.Text = "Converted " + $lq + "unknown" + $rq + " values in the regionType column to the " + $lq + "missing value" + $rq + " to correctly report on frequency.  Also changes " + $lq + "0" + $rq + " values in serviceArea column to " + $lq + "missing value" + $rq + ". " + [string][char]11

# ---------------------------------------------------------------------------
# Step 3: paragraph 16 ("Converted ..." originally) becomes the old
# "Used the frequency reports ..." text. Its trailing line break (the run
# holding it) is left alone by trimming the replaced range before it.
# ---------------------------------------------------------------------------
$p16 = $d.Paragraphs(16)
$r16 = $p16.Range
$r16.End = $r16.End - 2
$r16.Text = "Used the frequency reports to detect anomalies in the imported CSV data.  Then added statements to convert these anomalies and make the data consistent. e.g. some true/false values had " + $lsq + "t" + $rsq + " and " + $lsq + "f" + $rsq + " values these were converted to " + $lsq + "true" + $rsq + " and " + $lsq + "false" + $rsq + "."

# ---------------------------------------------------------------------------
# Step 4: paragraph 15 ("Used the frequency reports ..." originally) becomes
# the brand-new "linux command" text. Its trailing line break is likewise
# left in place.
# ---------------------------------------------------------------------------
$p15 = $d.Paragraphs(15)
$r15 = $p15.Range
$r15.End = $r15.End - 2
$r15.Text = "Used linux command (via git bash) " + $lq + "less" + $rq + " to view data any number of lines in, for fixing  import errors."

# ---------------------------------------------------------------------------
# Step 5: the final (previously empty, last list-numbered) paragraph is no
# longer part of the numbered list - strip numbering/list style and indent
# it under the list instead.
# ---------------------------------------------------------------------------
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$pLast.Range.ListFormat.RemoveNumbers()
$pLast.Style = "Normal"
$pLast.LeftIndent = 18
$pLast.Range.LanguageID = "en-IE"
